$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dCell = $ws.Cells.Item(2, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "60.699.38"
$dCell.Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  +2.51%  "

$dCell = $ws.Cells.Item(3, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "2.604.87"
$dCell.Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  +1.10%  "

$ws.Cells.Item(4, 5).Value = "  -0.01%  "

$dCell = $ws.Cells.Item(5, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "576.10"
$dCell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +3.78%  "

$dCell = $ws.Cells.Item(6, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "143.51"
$dCell.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +1.35%  "

$dCell = $ws.Cells.Item(7, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.996"
$dCell.Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  -0.26%  "

$ws.Cells.Item(8, 5).Value = "  +0.30%  "

$dCell = $ws.Cells.Item(9, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "2.629.39"
$dCell.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +1.84%  "

$dCell = $ws.Cells.Item(10, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "6.55"
$dCell.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  -1.82%  "

$ws.Cells.Item(11, 5).Value = "  +1.97%  "

$dCell = $ws.Cells.Item(12, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.157"
$dCell.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  -4.52%  "

$dCell = $ws.Cells.Item(13, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.369"
$dCell.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +5.07%  "

$dCell = $ws.Cells.Item(14, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "3.067.52"
$dCell.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +1.13%  "

$dCell = $ws.Cells.Item(15, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "60.660.05"
$dCell.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +2.51%  "

$dCell = $ws.Cells.Item(16, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "23.45"
$dCell.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +2.00%  "

$ws.Cells.Item(17, 5).Value = "  +4.22%  "

$dCell = $ws.Cells.Item(18, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "2.618.22"
$dCell.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +1.51%  "

$dCell = $ws.Cells.Item(19, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "11.31"
$dCell.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +9.49%  "

$dCell = $ws.Cells.Item(20, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "4.66"
$dCell.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +2.68%  "

$dCell = $ws.Cells.Item(21, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "349.81"
$dCell.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +3.70%  "

$ws.Cells.Item(22, 5).Value = "  +7.48%  "

$dCell = $ws.Cells.Item(23, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.998"
$dCell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -0.13%  "

$dCell = $ws.Cells.Item(24, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.519"
$dCell.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +10.45%  "

$dCell = $ws.Cells.Item(25, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "63.27"
$dCell.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +0.69%  "

$dCell = $ws.Cells.Item(26, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.996"
$dCell.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -0.40%  "

$dCell = $ws.Cells.Item(27, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.161"
$dCell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +1.34%  "

$dCell = $ws.Cells.Item(28, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "7.88"
$dCell.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +6.64%  "

$dCell = $ws.Cells.Item(29, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.0₃0797"
$dCell.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +3.40%  "

$ws.Cells.Item(30, 5).Value = "  +9.65%  "

$dCell = $ws.Cells.Item(32, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "6.36"
$dCell.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +3.27%  "

$dCell = $ws.Cells.Item(33, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "162.12"
$dCell.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +2.01%  "

$dCell = $ws.Cells.Item(34, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "19.57"
$dCell.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +2.62%  "

$dCell = $ws.Cells.Item(35, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "1.04"
$dCell.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +16.48%  "

$dCell = $ws.Cells.Item(36, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "4.26"
$dCell.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +5.77%  "

$dCell = $ws.Cells.Item(37, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "1.24"
$dCell.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  +6.67%  "

$dCell = $ws.Cells.Item(38, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "1.61"
$dCell.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +9.30%  "

$dCell = $ws.Cells.Item(39, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "37.95"
$dCell.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +1.84%  "

$dCell = $ws.Cells.Item(40, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "3.90"
$dCell.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +6.56%  "

$dCell = $ws.Cells.Item(41, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.850"
$dCell.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +0.04%  "

$dCell = $ws.Cells.Item(42, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "300.36"
$dCell.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +3.79%  "

$dCell = $ws.Cells.Item(43, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "134.25"
$dCell.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -1.70%  "

$ws.Cells.Item(44, 2).Value = "EnergySwap"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$dCell = $ws.Cells.Item(44, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "19.97"
$dCell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +5.78%  "

$ws.Cells.Item(45, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$dCell = $ws.Cells.Item(45, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.995"
$dCell.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -0.36%  "

$ws.Cells.Item(46, 2).Value = "RenderToken"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$dCell = $ws.Cells.Item(46, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "5.05"
$dCell.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +12.30%  "

$ws.Cells.Item(47, 2).Value = "Stellar"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$dCell = $ws.Cells.Item(47, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.0984"
$dCell.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +0.95%  "

$dCell = $ws.Cells.Item(48, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.607"
$dCell.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +2.71%  "

$ws.Cells.Item(49, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$dCell = $ws.Cells.Item(49, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "20.31"
$dCell.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +9.09%  "

$ws.Cells.Item(50, 2).Value = "Hedera"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$dCell = $ws.Cells.Item(50, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.0550"
$dCell.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +4.06%  "

$ws.Cells.Item(51, 2).Value = "VeChain"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$dCell = $ws.Cells.Item(51, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.0243"
$dCell.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +3.82%  "
